# Updated cryptos list on Fri Dec  8 02:43:30 UTC 2023 with GitHub Actions
#
# Refreshes price (col D) / 1h-volume-change (col E) figures for each coin
# row, and re-sorts a handful of rows whose rank order changed (Polygon /
# WrappedEther, BinanceUSD / FraxShare, FTXToken / TrustWalletToken, and
# Celestia replaced by TerraClassic) by swapping their Coin/Link/Price/
# Volume cells.
#
# Some "Price" strings look like plain numbers (e.g. "233.35", "0.460",
# "1.00") and Excel's COM Value setter auto-coerces those to numeric
# doubles, which would both change the stored cell type and silently drop
# meaningful trailing zeros (e.g. "0.460" -> 0.46, "1.00" -> 1). To keep
# them as literal text (matching the source data, which also has
# non-numeric "99.99.99"-style prices in the same column), force the
# range to Text format before assigning, then clear the formatting
# override afterwards so the cell's style reverts to the sheet default
# (avoids leaving a stray explicit "General"/Text number format behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.494.87'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '2.377.23'
$ws.Range("E3").Value = '  +5.74%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.79'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +10.48%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.460'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0962'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.41'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.27'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").Value = '2.730.39'
$ws.Range("E13").Value = '  +5.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("E16").Value = '  +1.43%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.431.14'
$ws.Range("E17").Value = '  +8.15%  '

$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.851'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("D19").Value = '43.508.54'
$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.09'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.77'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("E24").Value = '  +18.46%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("E27").Value = '  +1.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.90'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +9.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.94'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.59%  '

$ws.Range("E31").Value = '  +9.75%  '

$ws.Range("E32").Value = '  -8.83%  '

$ws.Range("E33").Value = '  +1.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0688'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.07'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.53%  '

$ws.Range("E37").Value = '  +2.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.44'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.06%  '

$ws.Range("E39").Value = '  -1.35%  '

$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("B41").Value = 'BinanceUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.92'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.99%  '

$ws.Range("E43").Value = '  +7.96%  '

$ws.Range("E44").Value = '  +10.78%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.53'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +4.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.72%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.22'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0948'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("D49").Value = '1.449.10'
$ws.Range("E49").Value = '  +0.88%  '

$ws.Range("D50").Value = '2.599.79'
$ws.Range("E50").Value = '  +5.81%  '

$ws.Range("B51").Value = 'TerraClassic'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000202'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -10.54%  '
